# University - Program - Import data : Start dates not exist in excel sheet
#
# Adds a new "Start Dates" column (AC) with per-program start-date lists,
# and fixes the "yes" flag for Program one's "Conditional Admissions"
# Ielts cell (K2), which had mistakenly been stored as the same string as
# the header-column YES constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the inserted "Start Dates" column
$ws.Range("AC1").Value2 = "Start Dates"

# New per-row values
$ws.Range("AC2").Value2 = "JAN,FEB,MAR,APR"
$ws.Range("AC3").Value2 = "JAN,FEB,MAR,APR,MAY,JUN,JUL,AUG"

# Program one's Ielts value corrected from "YES" to lowercase "yes"
$ws.Range("K2").Value2 = "yes"

# Give the new column a sensible custom width (matches the other
# description-style columns rather than the default 9.140625 width)
$ws.Columns("AC").ColumnWidth = 18.6

# Reflect the new selection position used while editing
$ws.Range("AB13").Select() | Out-Null
